$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 477, shifting existing rows 477:579 down to 478:580
$ws.Rows.Item(477).Insert()

# Populate the newly inserted row 477 with its data.
# Most columns repeat the values that were already present in the (now shifted)
# row 478 (the original row 477 content), except for D, M, N, O, P, Q, S, T
# which take new values.
$ws.Range("A477").Value = 10
$ws.Range("B477").Value = "Vega Modelo de Temuco"
$ws.Range("C477").Value = "La Araucanía"
$ws.Range("D477").Value = "10/21/2022"
$ws.Range("E477").Value = 9
$ws.Range("F477").Value = "Fruta"
$ws.Range("G477").Value = 100101
$ws.Range("H477").Value = "Berries"
$ws.Range("I477").Value = 100101007
$ws.Range("J477").Value = "Kiwi"
$ws.Range("K477").Value = "Hayward"
$ws.Range("L477").Value = "Primera"
$ws.Range("M477").Value = 285
$ws.Range("N477").Value = 10000
$ws.Range("O477").Value = 10000
$ws.Range("P477").Value = 10000
$ws.Range("Q477").Value = "$/bandeja 10 kilos"
$ws.Range("R477").Value = "Región de O'Higgins"
$ws.Range("S477").Value = 1000
$ws.Range("T477").Value = 10
